# Add a new weekly row (row 6) of data to the worksheet, mirroring the
# existing rows' structure/format (row 3 is the closest analogue: same
# Espárragos / "Sin especificar" variety, same $/kilo unit & origin).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

$ws.Cells.Item($row, 1).Value  = 7
$ws.Cells.Item($row, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value  = "Ñuble"

# Column D holds a date; match the existing date-formatted cells (style s="2",
# number format YYYY-MM-DD HH:MM:SS) by copying the format from the cell above
# and then setting the date value.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 4).Value  = "11/05/2021"

$ws.Cells.Item($row, 5).Value  = 16
$ws.Cells.Item($row, 6).Value  = 300000000
$ws.Cells.Item($row, 7).Value  = "Espárragos"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 440
$ws.Cells.Item($row, 11).Value = 900
$ws.Cells.Item($row, 12).Value = 1000
$ws.Cells.Item($row, 13).Value = 950
$ws.Cells.Item($row, 14).Value = "$/kilo"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 950
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
